# Apply the changes described by the commit:
#  - add a new Item row (Equip_Weapon_1 / 开山斧) to the Item.xlsx "表1" table
#  - grow the table/autofilter range to include the new row
#  - move the active selection to K13 (matches the post-edit cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- append the new data row (row 9) ------------------------------------
$ws.Range("A9").Value = "Equip_Weapon_1"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "开山斧"
$ws.Range("F9").Value = "开山斧武器"

# "Icon" column stores the id as text ("50004"), like the rest of column G
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "50004"

$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 10000
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 100

# --- grow the XML-mapped table (and its autofilter) to A1:K9 ------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K9"))

# --- match the saved cursor position -------------------------------------
$ws.Range("K13").Select()
